$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Total gain (%)" column (D),
# shifting it (and its formatting) to column E.
$ws.Range("D1").EntireColumn.Insert()

# Header for the new "Difference" column - match the bold/bordered/
# centered header style used by the other header cells.
$ws.Range("D1").Value = "Difference"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# Difference values (Current price - Purchase price), entered as text
# (leading apostrophe) so they match the source formatting exactly.
$diffs = @(
    "52.07",
    "19.12",
    "79.57",
    "45.34",
    "15.15",
    "1632.62",
    "87.00",
    "19.98",
    "28.17",
    "48.18",
    "-20.16",
    "2.76",
    "12.61",
    "38.05",
    "80.90",
    "27.45",
    "23.92",
    "9.97",
    "10.87",
    "-5.27",
    "62.03",
    "-6.06",
    "-4.85",
    "117.60",
    "19.29"
)

for ($i = 0; $i -lt $diffs.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = "'" + $diffs[$i]
    # Clear the implicit "quote prefix" style Excel applies when text is
    # entered with a leading apostrophe, so the cell stays unstyled like
    # its neighbors.
    $cell.Style = "Normal"
}
